$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (AC1) onto the new
# header cells so they match the bold/bordered header formatting, then set
# the header text for the new "Wins"/"Losses"/"Ties" columns (AD/AE/AF).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 77
    $ws.Cells.Item($r, 32).Value = 0
}
